# Add a "pid" column to the test data frame (Sheet1).
#
# The new column is inserted immediately before the existing column C
# ("m_wage_l1"), which pushes every column from C onward one slot to the
# right (C->D, D->E, ... Y->Z) and lets Excel rewrite all the dependent
# formulas automatically. A "pid" value (45-51) is then filled in for the
# seven data rows, and the active selection is moved to match the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new column; Excel shifts existing columns/formulas right.
$ws.Columns("C").Insert()

# Header for the newly inserted column.
$ws.Range("C1").Value = "pid"

# Sample "pid" values for the seven data rows.
$ws.Range("C2").Value = 45
$ws.Range("C3").Value = 46
$ws.Range("C4").Value = 47
$ws.Range("C5").Value = 48
$ws.Range("C6").Value = 49
$ws.Range("C7").Value = 50
$ws.Range("C8").Value = 51

# Move the active selection to where the author left off editing.
$ws.Range("E9").Select()
